$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Zoya", "Gomez", "zoya@gmail.com", "Madrid", "Spain", "all_excels_dir/sub_dir1/user_list4.xls"),
    @("Macie", "Lawson", "macie@gmail.com", "Paris", "France", "all_excels_dir/sub_dir1/user_list4.xls"),
    @("Sarah", "Oneal", "sarah@gmail.com", "Toronto", "Canada", "all_excels_dir/sub_dir1/user_list4.xls"),
    @("Carol", "Handley", "carol@gmail.com", "Vancouver", "Canada", "all_excels_dir/sub_dir1/user_list4.xls"),
    @("John", "Smith", "john@gmail.com", "New York", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Robert", "Williams", "bob@gmail.com", "London ", "UK", "all_excels_dir/user_list1.xlsx"),
    @("Sophia", "Miller", "sophia@gmail.com", "San Jose", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Suresh", "Govindaraj", "suresh@gmail.com", "Dallas", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Macy", "Barker", "macy@gmail.com", "New Jersey", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Casper", "Pitts", "casper@gmail.com", "Seattle", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Aleena", "Cobb", "aleena@gmail.com", "Chicago", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Murray", "Smart", "murray@gmail.com", "Phoenix", "USA", "all_excels_dir/user_list1.xlsx"),
    @("Raghu", "Srinivas", "raghu@gmail.com", "Chennai", "India", "all_excels_dir/user_list2.xls"),
    @("Emma", "Davis", "emma@gmail.com", "Seattle", "USA", "all_excels_dir/user_list2.xls"),
    @("Michael", "Brown", "michael@gmail.com", "Houston ", "USA", "all_excels_dir/user_list2.xls"),
    @("Yijiang", "Li", "yijiang@gmail.com", "Beijing", "China", "all_excels_dir/user_list2.xls"),
    @("Raj", "Kumar", "raj@gmail.com", "Mumbai", "India", "all_excels_dir/user_list3.csv"),
    @("Doug", "Wilson", "doug@gmail.com", "Chicago", "USA", "all_excels_dir/user_list3.csv"),
    @("Oliver", "Jones", "oliver@gmail.com", "Houston ", "USA", "all_excels_dir/user_list3.csv"),
    @("Kate", "Johnson", "kate@gmail.com", "Sydney", "Australia", "all_excels_dir/user_list3.csv"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
